$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) — update "想去人数" (interested count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 17
$ws1.Range("F5").Value = 6720
$ws1.Range("F9").Value = 6242
$ws1.Range("F12").Value = 1257
$ws1.Range("F16").Value = 123
$ws1.Range("F18").Value = 366
$ws1.Range("F21").Value = 4563
$ws1.Range("F22").Value = 60
$ws1.Range("F23").Value = 39
$ws1.Range("F24").Value = 25
$ws1.Range("F26").Value = 69

# Sheet "全部类型" (All types) — same events, shifted by one extra row (row 22 is a
# performance-type event not present in 展览), so rows 23-27 correspond to 展览's 22-26
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 17
$ws4.Range("F5").Value = 6720
$ws4.Range("F9").Value = 6242
$ws4.Range("F12").Value = 1257
$ws4.Range("F16").Value = 123
$ws4.Range("F18").Value = 366
$ws4.Range("F21").Value = 4563
$ws4.Range("F23").Value = 60
$ws4.Range("F24").Value = 39
$ws4.Range("F25").Value = 25
$ws4.Range("F27").Value = 69
